$wb = $excel.ActiveWorkbook
$st = $wb.Styles
Write-Output $st.Count
